$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected influent timing values (row 2) ---
$ws.Range("A2").Value = 0.4291666666666667
$ws.Range("B2").Value = 0.71666666673263535

# --- New column E: Rhodamine_meas (mg/L) ---
$ws.Range("E1").Value = "Rhodamine_meas (mg/L)"

$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 2.3548773174086999
$ws.Range("E4").Value = 0.78744553515667848
$ws.Range("E5").Value = 0.54142822447703309
$ws.Range("E6").Value = 0.44704147650757764
$ws.Range("E7").Value = 0.34435219052229044
$ws.Range("E8").Value = 0.33517570113637113
$ws.Range("E9").Value = 0.23336036652117151
$ws.Range("E10").Value = 0.26962934837980485
$ws.Range("E11").Value = 0.17917538157574339
$ws.Range("E12").Value = 0.14640220519746025
$ws.Range("E13").Value = 0.13460386170127833
$ws.Range("E14").Value = 0.12804922642562169
$ws.Range("E15").Value = 0.11537693155935219
$ws.Range("E16").Value = 0.10314161237812648
$ws.Range("E17").Value = 0.10051975826786386
$ws.Range("E18").Value = 0.11144415039395822
$ws.Range("E19").Value = 0.081292828125937744
$ws.Range("E20").Value = 0.076486095590456202
$ws.Range("E21").Value = 0.059007068188705193
$ws.Range("E22").Value = 0.053763359968179894
$ws.Range("E23").Value = 0.069931460314799576
$ws.Range("E24").Value = 0.058570092503661415
$ws.Range("E26").Value = 0.058570092503661415
$ws.Range("E27").Value = 0.05638521407844254
$ws.Range("E28").Value = 0.05201545722800479
$ws.Range("E29").Value = 0.05201545722800479
$ws.Range("E30").Value = 0.045023846267304393

# --- Column width / autofit for the new column ---
$ws.Columns.Item(5).AutoFit() | Out-Null

# --- Selection as left by the author ---
$ws.Range("C6").Select() | Out-Null
